{"js": "const body = context.document.body;\n\n// 1. Add a second line break right after the \"{{ formType }}\" merge field\n//    so the title block renders an extra blank line before the date fields.\nconst formTypeResults = body.search(\"{{ formType }}\", { matchCase: true });\nformTypeResults.load(\"items\");\nawait context.sync();\nformTypeResults.items[0].insertText(\"\\u000b\", Word.InsertLocation.end);\n\n// 2. Add \"modifiedAt\" and \"generatedAt\" merge fields (each on its own line)\n//    right after the existing \"createdAt\" merge field.\nconst createdAtResults = body.search(\"{{ createdAt }}\", { matchCase: true });\ncreatedAtResults.load(\"items\");\nawait context.sync();\ncreatedAtResults.items[0].insertText(\n  \"\\u000b{{ modifiedAt }}\\u000b{{ generatedAt }}\",\n  Word.InsertLocation.end\n);\nawait context.sync();\n\n// 3. The address/plots table's first two columns shift by 2 twips\n//    (3871/272 dxa -> 3869/274 dxa) across all rows.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst firstCell = table.getCell(0, 0);\nconst secondCell = table.getCell(0, 1);\nfirstCell.columnWidth = 193.45; // 3869 dxa\nsecondCell.columnWidth = 13.7; // 274 dxa\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Add a second line break right after the \"{{ formType }}\" merge field\n#    so the title block renders an extra blank line before the date fields.\n$range = $d.Content\n$range.Find.Execute(\"{{ formType }}\") | Out-Null\n$range.Text = $range.Text + [char]11\n\n# 2. Add \"modifiedAt\" and \"generatedAt\" merge fields (each on its own line)\n#    right after the existing \"createdAt\" merge field.\n$range2 = $d.Content\n$range2.Find.Execute(\"{{ createdAt }}\") | Out-Null\n$range2.Text = $range2.Text + [char]11 + \"{{ modifiedAt }}\" + [char]11 + \"{{ generatedAt }}\"\n\n# 3. The address/plots table's first two columns shift by 2 twips\n#    (3871/272 dxa -> 3869/274 dxa) across all rows.\n$t = $d.Tables.Item(1)\n$t.Columns.Item(1).Width = 193.45\n$t.Columns.Item(2).Width = 13.7\n"}
